# Insert two new data rows before the current row 198, pushing the existing
# rows 198-203 down to 200-205, and populate the two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 198 (existing 198..203 shift to 200..205)
$ws.Rows("198:199").Insert()

# New row 198
$ws.Cells.Item(198, 1).Value = 6
$ws.Cells.Item(198, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 44509
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = "Fruta"
$ws.Cells.Item(198, 7).Value = 100101
$ws.Cells.Item(198, 8).Value = "Berries"
$ws.Cells.Item(198, 9).Value = 100101001
$ws.Cells.Item(198, 10).Value = "Arándano (blue)"
$ws.Cells.Item(198, 11).Value = "Sin especificar"
$ws.Cells.Item(198, 12).Value = "Especial"
$ws.Cells.Item(198, 13).Value = 2000
$ws.Cells.Item(198, 14).Value = 7000
$ws.Cells.Item(198, 15).Value = 7000
$ws.Cells.Item(198, 16).Value = 7000
$ws.Cells.Item(198, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(198, 18).Value = "Provincia de Linares"
$ws.Cells.Item(198, 19).Value = 3500
$ws.Cells.Item(198, 20).Value = 2

# New row 199
$ws.Cells.Item(199, 1).Value = 6
$ws.Cells.Item(199, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(199, 3).Value = "Metropolitana"
$ws.Cells.Item(199, 4).Value = 44509
$ws.Cells.Item(199, 5).Value = 13
$ws.Cells.Item(199, 6).Value = "Fruta"
$ws.Cells.Item(199, 7).Value = 100101
$ws.Cells.Item(199, 8).Value = "Berries"
$ws.Cells.Item(199, 9).Value = 100101001
$ws.Cells.Item(199, 10).Value = "Arándano (blue)"
$ws.Cells.Item(199, 11).Value = "Sin especificar"
$ws.Cells.Item(199, 12).Value = "Segunda"
$ws.Cells.Item(199, 13).Value = 150
$ws.Cells.Item(199, 14).Value = 5000
$ws.Cells.Item(199, 15).Value = 5000
$ws.Cells.Item(199, 16).Value = 5000
$ws.Cells.Item(199, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(199, 18).Value = "Región Metropolitana"
$ws.Cells.Item(199, 19).Value = 2500
$ws.Cells.Item(199, 20).Value = 2

# Ensure the date column keeps the datetime number format used elsewhere
# in column D (style index 2 / numFmtId 165).
$ws.Range("D198:D199").NumberFormat = $ws.Range("D200").NumberFormat
